$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row2
$ws.Range("C2").Value = "Unknown Title"
$ws.Range("E2").Value = "[]"
$ws.Range("F2").Value = "not found"
$ws.Range("G2").Value = "N/A"
$ws.Range("H2").Value = "1970-01-01"
$ws.Range("J2").Value = ""

# row3
$ws.Range("D3").Value = "Uncovering risk factors playing roles in the severity of Coronavirus disease 2019 (Covid‐19) are important for understanding pathoimmunology of the disease caused by severe acute respiratory syndrome Coronavirus 2 (SARS CoV‐2).
 Genetic variations in innate immune genes have been found to be associated with Covid‐19 infections.
 A single‐nucleotide polymorphism (SNP) in a promoter region of tumor necrosis factor alpha (TNF‐α) gene, TNF‐α −308G&gt;A, increases expression of TNF‐α protein against infectious diseases leading to immune dysregulations and organ damage.
 This study aims to discover associations between TNF‐α −308G&gt;A SNP and Covid‐19 infection.
 Polymerase chain reaction‐restriction fragment length polymorphism (PCR‐RFLP) was used for genotyping a general Kurdish population and Covid‐19 patients.
 The homozygous mutant (AA) genotype was found to be rare in the current studied population.
 Interestingly, the heterozygous (GA) genotype was significantly (p value = 0.0342) higher in the Covid‐19 patients than the general population.
 This suggests that TNF‐α −308G&gt;A SNP might be associated with Covid‐19 infections.
 Further studies with larger sample sizes focusing on different ethnic populations are recommended.
"
$ws.Range("E3").Value = "[Hussein N.%Ali%NULL%0, Sherko S.%Niranji%sherko.subhan@garmian.edu.krd%1, Sirwan M. A.%Al‐Jaf%NULL%2, Sirwan M. A.%Al‐Jaf%NULL%0]"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "John Wiley and Sons Inc."

# row4
$ws.Range("D4").Value = "The COVID-19 pandemic has created an unprecedented need for epidemiological monitoring using diverse strategies.
 We conducted a project combining prevalence, seroprevalence, and genomic surveillance approaches to describe the initial pandemic stages in Betim City, Brazil.
 We collected 3239 subjects in a population-based age-, sex- and neighborhood-stratified, household, prospective; cross-sectional study divided into three surveys 21 days apart sampling the same geographical area.
 In the first survey, overall prevalence (participants positive in serological or molecular tests) reached 0.46% (90% CI 0.12–0.80%), followed by 2.69% (90% CI 1.88–3.49%) in the second survey and 6.67% (90% CI 5.42–7.92%) in the third.
 The underreporting reached 11, 19.6, and 20.4 times in each survey.
 We observed increased odds to test positive in females compared to males (OR 1.88 95% CI 1.25–2.82), while the single best predictor for positivity was ageusia/anosmia (OR 8.12, 95% CI 4.72–13.98).
 Thirty-five SARS-CoV-2 genomes were sequenced, of which 18 were classified as lineage B.
1.1.28, while 17 were B.
1.1.33. Multiple independent viral introductions were observed.
 Integration of multiple epidemiological strategies was able to adequately describe COVID-19 dispersion in the city.
 Presented results have helped local government authorities to guide pandemic management.
"
$ws.Range("E4").Value = "[Ana Valesca Fernandes Gilson%Silva%NULL%0, Diego%Menezes%NULL%1, Filipe Romero Rebello%Moreira%NULL%1, Octávio Alcântara%Torres%NULL%1, Paula Luize Camargos%Fonseca%NULL%1, Rennan Garcias%Moreira%NULL%1, Hugo José%Alves%NULL%1, Vivian Ribeiro%Alves%NULL%1, Tânia Maria de Resende%Amaral%NULL%1, Adriano Neves%Coelho%NULL%1, Júlia Maria%Saraiva Duarte%NULL%1, Augusto Viana%da Rocha%NULL%1, Luiz Gonzaga Paula%de Almeida%NULL%1, João Locke Ferreira%de Araújo%NULL%1, Hilton Soares%de Oliveira%NULL%1, Nova Jersey Cláudio%de Oliveira%NULL%1, Camila%Zolini%NULL%1, Jôsy Hubner%de Sousa%NULL%1, Elizângela Gonçalves%de Souza%NULL%1, Rafael Marques%de Souza%NULL%1, Luciana de Lima%Ferreira%NULL%1, Alexandra%Lehmkuhl Gerber%NULL%1, Ana Paula de Campos%Guimarães%NULL%1, Paulo Henrique Silva%Maia%NULL%1, Fernanda Martins%Marim%NULL%1, Lucyene%Miguita%NULL%1, Cristiane Campos%Monteiro%NULL%1, Tuffi Saliba%Neto%NULL%1, Fabrícia Soares Freire%Pugêdo%NULL%1, Daniel Costa%Queiroz%NULL%1, Damares Nigia Alborguetti Cuzzuol%Queiroz%NULL%1, Luciana Cunha%Resende-Moreira%NULL%1, Franciele Martins%Santos%NULL%1, Erika Fernanda Carlos%Souza%NULL%1, Carolina Moreira%Voloch%NULL%1, Ana Tereza%Vasconcelos%NULL%1, Renato Santana%de Aguiar%NULL%1, Renan Pedra%de Souza%NULL%1]"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "Frontiers Media S.A."
